$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$values = @(
    "47+48=",
    "90-67=",
    "32+37=",
    "20+9=",
    "74-50=",
    "86-5=",
    "34+49=",
    "38+39=",
    "29+19=",
    "82-77=",
    "26+59=",
    "2+96=",
    "80-3=",
    "51+3=",
    "36+33=",
    "87-41=",
    "6+32=",
    "29-8=",
    "66-30=",
    "28+49=",
    "77-53=",
    "90-32=",
    "70+26=",
    "5+2=",
    "23+6=",
    "43-11=",
    "66+30=",
    "88-55=",
    "31-22=",
    "99-79=",
    "8+56=",
    "75-65=",
    "85+1=",
    "81-76=",
    "25+16=",
    "32+47=",
    "83+10=",
    "76-24=",
    "51-47=",
    "56+41=",
    "67-13=",
    "1+82=",
    "73-19=",
    "9+24=",
    "7+12=",
    "8+63=",
    "65-59=",
    "57+35=",
    "24+12=",
    "60-18=",
    "9+14=",
    "28+43=",
    "17-4=",
    "21+70=",
    "17-1=",
    "82-44=",
    "39+13=",
    "23+26=",
    "10+42=",
    "64-41=",
    "81+7=",
    "13+5=",
    "59-1=",
    "65-57=",
    "27-21=",
    "78-30=",
    "92-44=",
    "70-35=",
    "66-23=",
    "7+57=",
    "88-83=",
    "91-38=",
    "58+39=",
    "65-21=",
    "65+6=",
    "80-39=",
    "16+20=",
    "3+86=",
    "23+49=",
    "15+2=",
    "13+25=",
    "86-64=",
    "30+52=",
    "88-71=",
    "31-28=",
    "43-43=",
    "45+27=",
    "90-57=",
    "56-12=",
    "89-51=",
    "36+48=",
    "94-63=",
    "71-40=",
    "6+77=",
    "9+54=",
    "28-21=",
    "75-53=",
    "64+9=",
    "5+70=",
    "41-33="
)

$idx = 0
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $row = $tbl.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "updated $idx cells"
